# Arduino tuner workbook edit:
# Insert a new "Pin ID" column before the existing "relay" column (old column C),
# shifting relay/component/value/wire# (old C:G) one column to the right (new D:H),
# and populate the new column C with the numeric pin-id values implied by the
# GPA/GPB pin naming (GPA0..GPA7 -> 0..7, GPB0..GPB7 -> 8..15) plus the Arduino
# digital-pin numbers for the two non-GPIO rows (D2 -> 2, D3 -> 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts old C:G to D:H.
$ws.Range("C1").EntireColumn.Insert()

# New column header.
$ws.Range("C1").Value = "Pin ID"

# New column C (Pin ID) values, keyed by row.
$pinIds = @{
    2  = 0
    3  = 8
    4  = 1
    5  = 9
    6  = 2
    7  = 10
    8  = 3
    10 = 11
    11 = 4
    12 = 12
    13 = 5
    14 = 13
    15 = 6
    16 = 14
    18 = 7
    19 = 15
    21 = 2
    22 = 3
}

foreach ($row in $pinIds.Keys) {
    $ws.Cells.Item($row, 3).Value = $pinIds[$row]
}

# Rows 24/25 have no data in any of columns C:H besides the "wire #" column
# (now G), so the freshly-inserted column C must stay blank there too;
# explicitly clear it so no stray empty-but-styled cell is left behind.
$ws.Cells.Item(24, 3).Clear()
$ws.Cells.Item(25, 3).Clear()

# Match the new column width seen in the edited file (stored width 15 in the
# OOXML character-width units), and update the selection to C22.
$ws.Columns.Item(3).ColumnWidth = 14.1
$ws.Range("C22").Select()
